# Insert a new "K.Barner" column before the existing "M.Evans" column (column J)
# on both the "Rushing" and "Receiving" sheets, shifting the later columns
# right by one. The new column gets the same header style as the other
# player columns and an "n" placeholder value in the data row, matching the
# existing (still-unfilled) columns.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $ws.Columns("J:J").Insert()
    $ws.Range("J1").Value = "K.Barner"
    $ws.Range("J2").Value = "n"
}
